$d = $word.ActiveDocument
$sec = $d.Sections(1)
$f1 = $sec.Footers(3)  # footer1.xml (even)
$fld = $f1.Range.Fields(2)
$fld.Delete()
Write-Output "deleted"
$p2 = $f1.Range.Paragraphs(2)
$r = $p2.Range
Write-Output ("p2 start=" + $r.Start + " end=" + $r.End)
$newField = $f1.Range.Fields.Add($r, 15, "TITLE  \* MERGEFORMAT", $false)
Write-Output ("newfield code=[" + $newField.Code.Text + "]")
